$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 19:07:57"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-26 19:07:52"
$wsZhCn.Range("K2").Value = "2016-08-26 19:08:17"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-26 19:08:24"
